$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to their rounded (2 decimal place) counterparts.
$ws.Range("B5").Value = 4.66
$ws.Range("C5").Value = 3.44
$ws.Range("D5").Value = 0.68
$ws.Range("E5").Value = 10.45
$ws.Range("F5").Value = 8.09
$ws.Range("H5").Value = 17.74
$ws.Range("J5").Value = 2.62
$ws.Range("K5").Value = 3.37
$ws.Range("L5").Value = 4.32
$ws.Range("M5").Value = 4.42
$ws.Range("N5").Value = 1.4
$ws.Range("Q5").Value = 3.42
$ws.Range("R5").Value = 0.57
$ws.Range("W5").Value = 6.94
$ws.Range("X5").Value = 3.68
$ws.Range("Z5").Value = 8.73
$ws.Range("AA5").Value = 3.07
$ws.Range("AB5").Value = 2.83
$ws.Range("AC5").Value = 3.31
$ws.Range("AD5").Value = 4.4
$ws.Range("AF5").Value = 16.5
$ws.Range("AG5").Value = 1.85

# Delete row 6 entirely (last row of data).
$ws.Rows("6:6").Delete()
